$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.609.47"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.853.27"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "263.95"
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5271"
$ws.Range("E7").Value = "  +1.90%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3243"
$ws.Range("E8").Value = "  +0.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06801"
$ws.Range("E9").Value = "  +1.09%  "
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7844"
$ws.Range("E11").Value = "  +2.49%  "
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.868.65"
$ws.Range("E13").Value = "  +2.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.74"
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.043"
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.00"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007958"
$ws.Range("E19").Value = "  +1.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.632.60"
$ws.Range("E20").Value = "  +1.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.642"
$ws.Range("E21").Value = "  +2.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.480"
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.019"
$ws.Range("E23").Value = "  +2.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.78"
$ws.Range("E24").Value = "  -0.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.171"
$ws.Range("E25").Value = "  -5.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.680"
$ws.Range("E26").Value = "  +1.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.02"
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "111.70"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.186"
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08724"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.105"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04865"
$ws.Range("E32").Value = "  +0.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7225"
$ws.Range("E33").Value = "  +6.47%  "
$ws.Range("E34").Value = "  +0.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.873"
$ws.Range("E35").Value = "  +0.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.113"
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("E37").Value = "  +3.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01795"
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.4878"
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9014"
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "111.71"
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.970"
$ws.Range("E42").Value = "  -2.66%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4195"
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05886"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.015"
$ws.Range("E47").Value = "  -0.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1239"
$ws.Range("E48").Value = "  -1.26%  "
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.8908"
$ws.Range("E50").Value = "  +3.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "60.06"
$ws.Range("E51").Value = "  +1.62%  "
